# Automatische test-sync: 2025-08-05 17:21:50
# Append a new "Logs" row (row 17) for the testmail "Kun jij dit even regelen?"
# and bump the matching "Planning / Afspraak" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 17 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "Kun jij dit even regelen?"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D17").Value = "Planning / Afspraak"
$logs.Range("E17").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F17").Value = "2025-08-05 17:21:18"
$logs.Range("G17").Value = "Ja"
$logs.Range("H17").Value = "Ja"
$logs.Range("I17").Value = "Nee"
$logs.Range("J17").Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J 2:16 -> 2:17) so the new
# row is covered, same way Excel grows them when a table/range is extended.
$logs.Range("D2:D16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D17"))
$logs.Range("G2:G16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G17"))
$logs.Range("H2:H16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H17"))
$logs.Range("I2:I16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I17"))
$logs.Range("J2:J16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J17"))

# --- Sheet "Dashboard": bump the "Planning / Afspraak" count ---------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 11
